$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'42.726.84"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'  +1.77%  "
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'2.313.40"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'  +1.02%  "
$ws.Range("E3").Style = "Normal"
$ws.Range("E4").Value = "'  +0.10%  "
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'317.28"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'  -0.22%  "
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = "'105.02"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'  +1.42%  "
$ws.Range("E6").Style = "Normal"
$ws.Range("D7").Value = "'0.630"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "'  +0.16%  "
$ws.Range("E7").Style = "Normal"
$ws.Range("E8").Value = "'  +0.21%  "
$ws.Range("E8").Style = "Normal"
$ws.Range("D9").Value = "'0.610"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'  +0.74%  "
$ws.Range("E9").Style = "Normal"
$ws.Range("D10").Value = "'40.21"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'  +2.20%  "
$ws.Range("E10").Style = "Normal"
$ws.Range("D11").Value = "'0.0909"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'  +0.55%  "
$ws.Range("E11").Style = "Normal"
$ws.Range("E12").Value = "'  +3.68%  "
$ws.Range("E12").Style = "Normal"
$ws.Range("E13").Value = "'  +0.98%  "
$ws.Range("E13").Style = "Normal"
$ws.Range("D14").Value = "'0.987"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'  +2.79%  "
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Value = "'15.47"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'  +1.61%  "
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").Value = "'2.667.06"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'  +1.14%  "
$ws.Range("E16").Style = "Normal"
$ws.Range("D17").Value = "'2.317.60"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'  +1.08%  "
$ws.Range("E17").Style = "Normal"
$ws.Range("D18").Value = "'42.669.26"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'  +1.48%  "
$ws.Range("E18").Style = "Normal"
$ws.Range("D19").Value = "'7.66"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'  +3.30%  "
$ws.Range("E19").Style = "Normal"
$ws.Range("D20").Value = "'0.0000106"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'  +0.85%  "
$ws.Range("E20").Style = "Normal"
$ws.Range("D21").Value = "'13.60"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'  +35.43%  "
$ws.Range("E21").Style = "Normal"
$ws.Range("D22").Value = "'73.99"
$ws.Range("D22").Style = "Normal"
$ws.Range("D23").Value = "'3.54"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'  -2.26%  "
$ws.Range("E23").Style = "Normal"
$ws.Range("D24").Value = "'269.71"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'  -3.51%  "
$ws.Range("E24").Style = "Normal"
$ws.Range("D25").Value = "'2.24"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'  -0.64%  "
$ws.Range("E25").Style = "Normal"
$ws.Range("E26").Value = "'  -0.31%  "
$ws.Range("E26").Style = "Normal"
$ws.Range("D27").Value = "'10.94"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "'  +1.40%  "
$ws.Range("E27").Style = "Normal"
$ws.Range("D28").Value = "'2.34"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "'  -1.99%  "
$ws.Range("E28").Style = "Normal"
$ws.Range("D29").Value = "'22.68"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "'  -1.04%  "
$ws.Range("E29").Style = "Normal"
$ws.Range("D30").Value = "'38.39"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "'  +6.27%  "
$ws.Range("E30").Style = "Normal"
$ws.Range("D31").Value = "'6.59"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "'  +13.27%  "
$ws.Range("E31").Style = "Normal"
$ws.Range("D32").Value = "'166.47"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "'  +2.12%  "
$ws.Range("E32").Style = "Normal"
$ws.Range("D33").Value = "'0.0888"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "'  +2.03%  "
$ws.Range("E33").Style = "Normal"
$ws.Range("B34").Value = "'WEMIXToken"
$ws.Range("B34").Style = "Normal"
$ws.Range("C34").Value = "'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("C34").Style = "Normal"
$ws.Range("D34").Value = "'2.67"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "'  -6.34%  "
$ws.Range("E34").Style = "Normal"
$ws.Range("B35").Value = "'Stellar"
$ws.Range("B35").Style = "Normal"
$ws.Range("C35").Value = "'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("C35").Style = "Normal"
$ws.Range("D35").Value = "'0.132"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "'  -2.92%  "
$ws.Range("E35").Style = "Normal"
$ws.Range("D36").Value = "'0.114"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "'  +0.39%  "
$ws.Range("E36").Style = "Normal"
$ws.Range("D37").Value = "'4.60"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "'  +1.97%  "
$ws.Range("E37").Style = "Normal"
$ws.Range("D38").Value = "'0.0354"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "'  +1.50%  "
$ws.Range("E38").Style = "Normal"
$ws.Range("B39").Value = "'LidoDAOToken"
$ws.Range("B39").Style = "Normal"
$ws.Range("C39").Value = "'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("C39").Style = "Normal"
$ws.Range("D39").Value = "'2.80"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "'  -3.88%  "
$ws.Range("E39").Style = "Normal"
$ws.Range("B40").Value = "'NEARProtocol"
$ws.Range("B40").Style = "Normal"
$ws.Range("C40").Value = "'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("C40").Style = "Normal"
$ws.Range("D40").Value = "'3.72"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'  +0.10%  "
$ws.Range("E40").Style = "Normal"
$ws.Range("D41").Value = "'1.65"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'  +13.36%  "
$ws.Range("E41").Style = "Normal"
$ws.Range("D42").Value = "'99.90"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'  +0.70%  "
$ws.Range("E42").Style = "Normal"
$ws.Range("D43").Value = "'70.54"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'  +1.62%  "
$ws.Range("E43").Style = "Normal"
$ws.Range("D44").Value = "'0.226"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'  +0.65%  "
$ws.Range("E44").Style = "Normal"
$ws.Range("E45").Value = "'  +0.02%  "
$ws.Range("E45").Style = "Normal"
$ws.Range("D46").Value = "'117.57"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'  +4.15%  "
$ws.Range("E46").Style = "Normal"
$ws.Range("D47").Value = "'12.33"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'  +3.94%  "
$ws.Range("E47").Style = "Normal"
$ws.Range("D48").Value = "'81.97"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'  +6.81%  "
$ws.Range("E48").Style = "Normal"
$ws.Range("D49").Value = "'1.650.65"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "'  +4.57%  "
$ws.Range("E49").Style = "Normal"
$ws.Range("D50").Value = "'5.31"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "'  +0.67%  "
$ws.Range("E50").Style = "Normal"
$ws.Range("E51").Value = "'  -0.91%  "
$ws.Range("E51").Style = "Normal"
